# Y4_B2526_Excuses.xlsx - attendance app upload
#
# The sheet is a flat log: column A holds a "Student ID" that LOOKS numeric
# but is actually stored as text (the source app writes every cell as a
# plain string). This script:
#   1. Refreshes the Student ID in column A for the existing rows (2-26)
#      with newer log IDs.
#   2. Appends six brand-new log rows (27-32) with the same
#      Subject/Date/Time/Type/User values the app always stamps.
#
# Because Excel auto-detects a purely-numeric literal typed into a cell and
# stores it as a Number, we quote-prefix the literal (leading "'") to force
# Text storage like the source file expects. Typing a quote-prefixed literal
# also flips the cell onto an auxiliary "quote-prefixed" style, so
# immediately after, we re-apply the row's correct alternating banding style
# by pasting formats from the matching cell in column B of the same row
# (column B is never touched by this script, so it always holds the
# untouched/original style for that row).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- 1. Update existing Student IDs (column A, rows 2-26) ---------------
$studentIdUpdates = @(
    @{ Row = 2;  Id = "201065" }
    @{ Row = 3;  Id = "201080" }
    @{ Row = 4;  Id = "201157" }
    @{ Row = 5;  Id = "201171" }
    @{ Row = 6;  Id = "201190" }
    @{ Row = 7;  Id = "201197" }
    @{ Row = 8;  Id = "201218" }
    @{ Row = 9;  Id = "201237" }
    @{ Row = 10; Id = "201252" }
    @{ Row = 11; Id = "201253" }
    @{ Row = 12; Id = "201255" }
    @{ Row = 13; Id = "201297" }
    @{ Row = 14; Id = "201328" }
    @{ Row = 15; Id = "201337" }
    @{ Row = 16; Id = "201397" }
    @{ Row = 17; Id = "201398" }
    @{ Row = 18; Id = "201438" }
    @{ Row = 19; Id = "201465" }
    @{ Row = 20; Id = "201495" }
    @{ Row = 21; Id = "201501" }
    @{ Row = 22; Id = "201513" }
    @{ Row = 23; Id = "201529" }
    @{ Row = 24; Id = "201560" }
    @{ Row = 25; Id = "201563" }
    @{ Row = 26; Id = "201564" }
)

foreach ($u in $studentIdUpdates) {
    $r = $u.Row

    $cell = $ws.Cells.Item($r, 1)
    $cell.Value = "'" + $u.Id

    # Restore the row's real style (quote-prefixing a literal switches the
    # cell onto a separate style); column B of the same row always keeps the
    # original, untouched banding style for that row.
    $ws.Cells.Item($r, 2).Copy() | Out-Null
    $cell.PasteSpecial(-4122) | Out-Null  # xlPasteFormats
}

# --- 2. Append six new log rows (27-32) ----------------------------------
$newRows = @(
    @{ Row = 27; Id = "201572" }
    @{ Row = 28; Id = "201574" }
    @{ Row = 29; Id = "201632" }
    @{ Row = 30; Id = "201638" }
    @{ Row = 31; Id = "201669" }
    @{ Row = 32; Id = "201670" }
)

foreach ($entry in $newRows) {
    $r = $entry.Row

    # New rows start with the default "General" style (no banding) — paint
    # them with the correct alternating banding style *before* typing any
    # values, by copying the format from the nearest existing row two rows
    # up (same parity, e.g. row 25 -> row 27, row 26 -> row 28, etc.), which
    # mirrors the existing banding pattern (even row -> style 2, odd row ->
    # style 3).
    $styleSourceRow = $r - 2
    $ws.Range("A" + $styleSourceRow + ":F" + $styleSourceRow).Copy() | Out-Null
    $ws.Range("A" + $r + ":F" + $r).PasteSpecial(-4122) | Out-Null  # xlPasteFormats

    # Same fixed values the app stamps on every one of these new rows.
    # Date/time-looking literals are quote-prefixed so Excel stores them as
    # plain text (matching the source file) instead of auto-converting them
    # to date/time serial numbers.
    $ws.Cells.Item($r, 2).Value = "general surgery"
    $ws.Cells.Item($r, 3).Value = "'08/09/2025"
    $ws.Cells.Item($r, 4).Value = "'10:30:00"
    $ws.Cells.Item($r, 5).Value = "Excuse"
    $ws.Cells.Item($r, 6).Value = "System"

    # Column A (Student ID) as text, same as the existing rows.
    $ws.Cells.Item($r, 1).Value = "'" + $entry.Id

    # Quote-prefixing nudges a cell onto an auxiliary style, so re-apply the
    # row's banding style once more now that every value is in place.
    $ws.Range("A" + $styleSourceRow + ":F" + $styleSourceRow).Copy() | Out-Null
    $ws.Range("A" + $r + ":F" + $r).PasteSpecial(-4122) | Out-Null  # xlPasteFormats
}

$excel.CutCopyMode = $false
